# Case 4_6 (380 kV case) results: update res_bus/vm_pu.xlsx values.
# Column B is the slack/reference voltage set-point (dropped from 1.05 pu to
# 1.02 pu); columns C-F and I-N are the recalculated bus voltage magnitudes
# (column H has no bus and stays empty) for time steps in rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021703328107251
$ws.Range("D2").Value = 1.025682068568814
$ws.Range("E2").Value = 1.025323746927289
$ws.Range("F2").Value = 1.020138785652877
$ws.Range("I2").Value = 1.027792440572788
$ws.Range("J2").Value = 1.026893301106865
$ws.Range("K2").Value = 1.028506925652289
$ws.Range("L2").Value = 1.028149652795205
$ws.Range("M2").Value = 1.022979957318763
$ws.Range("N2").Value = 1.01288847942797

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023158140410001
$ws.Range("D3").Value = 1.026715445042946
$ws.Range("E3").Value = 1.026726601252953
$ws.Range("F3").Value = 1.022240355144684
$ws.Range("I3").Value = 1.028046694830089
$ws.Range("J3").Value = 1.027983218889189
$ws.Range("K3").Value = 1.029347043456439
$ws.Range("L3").Value = 1.02935816940541
$ws.Range("M3").Value = 1.024884149848731
$ws.Range("N3").Value = 1.013261333857136

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024097285015622
$ws.Range("D4").Value = 1.02738210194127
$ws.Range("E4").Value = 1.027632529772037
$ws.Range("F4").Value = 1.023597449895691
$ws.Range("I4").Value = 1.028209022641972
$ws.Range("J4").Value = 1.028685910111212
$ws.Range("K4").Value = 1.029888057574872
$ws.Range("L4").Value = 1.030137840757873
$ws.Range("M4").Value = 1.026113189372571
$ws.Range("N4").Value = 1.013501418083647

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024491580937568
$ws.Range("D5").Value = 1.027661890065076
$ws.Range("E5").Value = 1.028012957350424
$ws.Range("F5").Value = 1.024167331847763
$ws.Range("I5").Value = 1.02827674265516
$ws.Range("J5").Value = 1.028980716818536
$ws.Range("K5").Value = 1.030114883612431
$ws.Range("L5").Value = 1.03046506755829
$ws.Range("M5").Value = 1.026629153822062
$ws.Range("N5").Value = 1.013602070110541

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02455775472906
$ws.Range("D6").Value = 1.027708840073896
$ws.Range("E6").Value = 1.028076808166975
$ws.Range("F6").Value = 1.024262980553591
$ws.Range("I6").Value = 1.028288082558359
$ws.Range("J6").Value = 1.029030180928823
$ws.Range("K6").Value = 1.030152932727725
$ws.Range("L6").Value = 1.030519978539362
$ws.Range("M6").Value = 1.026715744615017
$ws.Range("N6").Value = 1.013618953718345

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024102555649716
$ws.Range("D7").Value = 1.027385842340759
$ws.Range("E7").Value = 1.027637614721928
$ws.Range("F7").Value = 1.023605067178698
$ws.Range("I7").Value = 1.028209929571007
$ws.Range("J7").Value = 1.028689851699278
$ws.Range("K7").Value = 1.029891090850111
$ws.Range("L7").Value = 1.030142215317694
$ws.Range("M7").Value = 1.026120086536078
$ws.Range("N7").Value = 1.013502764095653

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022195453727851
$ws.Range("D8").Value = 1.026031722128416
$ws.Range("E8").Value = 1.025798228770371
$ws.Range("F8").Value = 1.020849604943952
$ws.Range("I8").Value = 1.027878822063897
$ws.Range("J8").Value = 1.027262178399271
$ws.Range("K8").Value = 1.028791388857837
$ws.Range("L8").Value = 1.028558561974371
$ws.Range("M8").Value = 1.023624141216134
$ws.Range("N8").Value = 1.013014732665672

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018817467339376
$ws.Range("D9").Value = 1.023629948871508
$ws.Range("E9").Value = 1.022542710133389
$ws.Range("F9").Value = 1.015972062995973
$ws.Range("I9").Value = 1.027278487537563
$ws.Range("J9").Value = 1.024726487939678
$ws.Range("K9").Value = 1.026833414510196
$ws.Range("L9").Value = 1.025749821446003
$ws.Range("M9").Value = 1.019201389178886
$ws.Range("N9").Value = 1.012145626630845

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016553117362615
$ws.Range("D10").Value = 1.022017868144466
$ws.Range("E10").Value = 1.020362185706901
$ws.Range("F10").Value = 1.01270424496136
$ws.Range("I10").Value = 1.026866780036499
$ws.Range("J10").Value = 1.023022118740636
$ws.Range("K10").Value = 1.025514178519383
$ws.Range("L10").Value = 1.023864602832196
$ws.Range("M10").Value = 1.016235194897558
$ws.Range("N10").Value = 1.011559921669934

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0155695537369
$ws.Range("D11").Value = 1.0213171480227
$ws.Range("E11").Value = 1.019415451458833
$ws.Range("F11").Value = 1.011285131537907
$ws.Range("I11").Value = 1.026685751600657
$ws.Range("J11").Value = 1.022280697335263
$ws.Range("K11").Value = 1.024939550683021
$ws.Range("L11").Value = 1.023045145576933
$ws.Range("M11").Value = 1.014946342227964
$ws.Range("N11").Value = 1.011304774168219

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.015203738418977
$ws.Range("D12").Value = 1.021056459390797
$ws.Range("E12").Value = 1.019063397469426
$ws.Range("F12").Value = 1.010757364755753
$ws.Range("I12").Value = 1.026618092758377
$ws.Range("J12").Value = 1.022004777117766
$ws.Range("K12").Value = 1.024725591797454
$ws.Range("L12").Value = 1.022740279813723
$ws.Range("M12").Value = 1.014466910643693
$ws.Range("N12").Value = 1.011209767473225

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.015282228788167
$ws.Range("D13").Value = 1.0211123966621
$ws.Range("E13").Value = 1.019138932289241
$ws.Range("F13").Value = 1.010870602106529
$ws.Range("I13").Value = 1.026632624715598
$ws.Range("J13").Value = 1.022063986778823
$ws.Range("K13").Value = 1.0247715101953
$ws.Range("L13").Value = 1.022805696572021
$ws.Range("M13").Value = 1.014569782131006
$ws.Range("N13").Value = 1.011230157352004

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015539325111479
$ws.Range("D14").Value = 1.0212956078266
$ws.Range("E14").Value = 1.019386358684063
$ws.Range("F14").Value = 1.011241519471267
$ws.Range("I14").Value = 1.026680167418462
$ws.Range("J14").Value = 1.022257900415823
$ws.Range("K14").Value = 1.024921875353783
$ws.Range("L14").Value = 1.02301995518807
$ws.Range("M14").Value = 1.014906726559227
$ws.Range("N14").Value = 1.011296925668103

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015697667191708
$ws.Range("D15").Value = 1.021408435717716
$ws.Range("E15").Value = 1.019538753608022
$ws.Range("F15").Value = 1.011469967700215
$ws.Range("I15").Value = 1.02670940473236
$ws.Range("J15").Value = 1.022377307372519
$ws.Range("K15").Value = 1.025014451628861
$ws.Range("L15").Value = 1.023151902699148
$ws.Range("M15").Value = 1.015114236457489
$ws.Range("N15").Value = 1.011338032787545

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016618327178187
$ws.Range("D16").Value = 1.022064315554549
$ws.Range("E16").Value = 1.020424962602938
$ws.Range("F16").Value = 1.01279833778416
$ws.Range("I16").Value = 1.026878735993767
$ws.Range("J16").Value = 1.023071251589614
$ws.Range("K16").Value = 1.025552242639862
$ws.Range("L16").Value = 1.023918920367817
$ws.Range("M16").Value = 1.016320635683025
$ws.Range("N16").Value = 1.011576822392238

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.01719499811278
$ws.Range("D17").Value = 1.022475009222619
$ws.Range("E17").Value = 1.02098016676685
$ws.Range("F17").Value = 1.013630466943759
$ws.Range("I17").Value = 1.026984213161727
$ws.Range("J17").Value = 1.02350562245178
$ws.Range("K17").Value = 1.025888672078583
$ws.Range("L17").Value = 1.024399200986319
$ws.Range("M17").Value = 1.017076165164919
$ws.Range("N17").Value = 1.011726196039876

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01753106369175
$ws.Range("D18").Value = 1.022714301967889
$ws.Range("E18").Value = 1.02130376242016
$ws.Range("F18").Value = 1.014115437337601
$ws.Range("I18").Value = 1.027045470444488
$ws.Range("J18").Value = 1.023758654474706
$ws.Range("K18").Value = 1.026084579090409
$ws.Range("L18").Value = 1.024679037787807
$ws.Range("M18").Value = 1.017516422885103
$ws.Range("N18").Value = 1.011813175428551

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017645603399015
$ws.Range("D19").Value = 1.022795851088976
$ws.Range("E19").Value = 1.021414058864249
$ws.Range("F19").Value = 1.014280733210276
$ws.Range("I19").Value = 1.0270663125978
$ws.Range("J19").Value = 1.023844876318331
$ws.Range("K19").Value = 1.026151323140681
$ws.Range("L19").Value = 1.024774403863011
$ws.Range("M19").Value = 1.017666467161364
$ws.Range("N19").Value = 1.011842808189749

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017133157570807
$ws.Range("D20").Value = 1.022430972411612
$ws.Range("E20").Value = 1.020920624093024
$ws.Range("F20").Value = 1.013541228547996
$ws.Range("I20").Value = 1.026972923968692
$ws.Range("J20").Value = 1.023459052685166
$ws.Range("K20").Value = 1.02585261021642
$ws.Range("L20").Value = 1.024347702811525
$ws.Range("M20").Value = 1.016995148625586
$ws.Range("N20").Value = 1.011710184964797

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015463629892966
$ws.Range("D21").Value = 1.021241668083703
$ws.Range("E21").Value = 1.019313508765056
$ws.Range("F21").Value = 1.011132311555151
$ws.Range("I21").Value = 1.026666178805662
$ws.Range("J21").Value = 1.022200812207182
$ws.Range("K21").Value = 1.024877610879333
$ws.Range("L21").Value = 1.022956874778536
$ws.Range("M21").Value = 1.014807524172627
$ws.Range("N21").Value = 1.011277270536083

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014411171335482
$ws.Range("D22").Value = 1.020491529420833
$ws.Range("E22").Value = 1.018300762101703
$ws.Range("F22").Value = 1.009613986265582
$ws.Range("I22").Value = 1.026470902981808
$ws.Range("J22").Value = 1.021406673944395
$ws.Range("K22").Value = 1.024261598019425
$ws.Range("L22").Value = 1.022079607153141
$ws.Range("M22").Value = 1.01342804935877
$ws.Range("N22").Value = 1.011003727186289

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.014969365374107
$ws.Range("D23").Value = 1.020889419869259
$ws.Range("E23").Value = 1.018837859061156
$ws.Range("F23").Value = 1.01041924246229
$ws.Range("I23").Value = 1.026574652025284
$ws.Range("J23").Value = 1.021827952392576
$ws.Range("K23").Value = 1.024588444140407
$ws.Range("L23").Value = 1.022544932097053
$ws.Range("M23").Value = 1.014159724528977
$ws.Range("N23").Value = 1.011148866999202

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017161101566969
$ws.Range("D24").Value = 1.022450871548036
$ws.Range("E24").Value = 1.020947529623083
$ws.Range("F24").Value = 1.01358155276231
$ws.Range("I24").Value = 1.02697802589032
$ws.Range("J24").Value = 1.023480096575183
$ws.Range("K24").Value = 1.02586890602885
$ws.Range("L24").Value = 1.024370973555513
$ws.Range("M24").Value = 1.017031757840095
$ws.Range("N24").Value = 1.011717420136946

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019692888770926
$ws.Range("D25").Value = 1.024252757413305
$ws.Range("E25").Value = 1.023386091337063
$ws.Range("F25").Value = 1.017235771092883
$ws.Range("I25").Value = 1.02743570262522
$ws.Range("J25").Value = 1.025384439443328
$ws.Range("K25").Value = 1.027342024152995
$ws.Range("L25").Value = 1.02647815198422
$ws.Range("M25").Value = 1.020347808665140
$ws.Range("N25").Value = 1.012371410037132
